$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("I4").Value = 0
$ws1.Range("L4").Value = 0
$ws1.Range("P4").Value = 0
$ws1.Range("L5").Value = 0
$ws1.Range("M5").Value = 0
$ws1.Range("E6").Value = 0
$ws1.Range("M6").Value = 0
$ws1.Range("M12").Value = 0
$ws1.Range("C15").Value = 0
$ws1.Range("D15").Value = 0
$ws1.Range("L15").Value = 0
$ws1.Range("M15").Value = 0
$ws1.Range("P15").Value = 0
$ws1.Range("L26").Value = 0
$ws1.Range("H27").Value = 0
$ws1.Range("K27").Value = 0
$ws1.Range("M27").Value = 0
$ws1.Range("D28").Value = 0
$ws1.Range("K28").Value = 0
$ws1.Range("L28").Value = 0
$ws1.Range("M28").Value = 0
$ws1.Range("P28").Value = 0
$ws1.Range("M29").Value = 0
$ws1.Range("C30").Value = 0
$ws1.Range("D30").Value = 0
$ws1.Range("I31").Value = 0
$ws1.Range("K31").Value = 0
$ws1.Range("L31").Value = 0
$ws1.Range("M31").Value = 0
$ws1.Range("D32").Value = 0
$ws1.Range("D34").Value = 0
$ws1.Range("K34").Value = 0
$ws1.Range("M34").Value = 0
$ws1.Range("P34").Value = 0
$ws1.Range("L40").Value = 0
$ws1.Range("C41").Value = 0
$ws1.Range("D41").Value = 0
$ws1.Range("L41").Value = 0
$ws1.Range("M41").Value = 0
$ws1.Range("D47").Value = 0
$ws1.Range("E47").Value = 0
$ws1.Range("L47").Value = 0
$ws1.Range("D48").Value = 0
$ws1.Range("H48").Value = 0
$ws1.Range("I48").Value = 0
$ws1.Range("M48").Value = 0
$ws1.Range("N48").Value = 0
$ws1.Range("I49").Value = 0
$ws1.Range("L49").Value = 0
$ws1.Range("M49").Value = 0
$ws1.Range("D50").Value = 0
$ws1.Range("H50").Value = 0
$ws1.Range("H52").Value = 0
$ws1.Range("L52").Value = 0
$ws1.Range("E56").Value = 0
$ws1.Range("P56").Value = 0
$ws1.Range("I58").Value = 0
$ws1.Range("K58").Value = 0
$ws1.Range("C60").Value = "0 de 58"
$ws1.Range("D60").Value = "0 de 58"
$ws1.Range("E60").Value = "0 de 58"
$ws1.Range("H60").Value = "0 de 58"
$ws1.Range("I60").Value = "0 de 58"
$ws1.Range("K60").Value = "0 de 58"
$ws1.Range("L60").Value = "0 de 58"
$ws1.Range("M60").Value = "0 de 58"
$ws1.Range("P60").Value = "0 de 58"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item(2)

# Column width adjustments (ColumnWidth setter adds ~0.8333 padding vs stored width, so compensate)
$ws2.Columns.Item(4).ColumnWidth = 16 - 0.8333333333333333
$ws2.Columns.Item(5).ColumnWidth = 14 - 0.8333333333333333
$ws2.Columns.Item(6).ColumnWidth = 15 - 0.8333333333333333

$ws2.Range("C1").Value = "agosto"
$ws2.Range("D1").Value = "septiembre"
$ws2.Range("E1").Value = "octubre"
$ws2.Range("F1").Value = "noviembre"
$ws2.Range("D4").Value = 794.99
$ws2.Range("E4").Value = 1405.49
$ws2.Range("F4").Value = 0
$ws2.Range("C5").Value = 5652.48
$ws2.Range("D5").Value = 0
$ws2.Range("E5").Value = 21118.81
$ws2.Range("F5").Value = 0
$ws2.Range("C6").Value = 1944.78
$ws2.Range("D6").Value = 1603.38
$ws2.Range("E6").Value = 387.91
$ws2.Range("F6").Value = 0
$ws2.Range("C12").Value = -295.8
$ws2.Range("D12").Value = 0
$ws2.Range("E12").Value = 3404.07
$ws2.Range("F12").Value = 0
$ws2.Range("C13").Value = 0
$ws2.Range("C15").Value = 11261.15
$ws2.Range("D15").Value = 8223.309999999999
$ws2.Range("E15").Value = 8398.709999999999
$ws2.Range("F15").Value = 0
$ws2.Range("C16").Value = 174.18
$ws2.Range("D16").Value = 0
$ws2.Range("C24").Value = 8385.889999999999
$ws2.Range("D24").Value = 146.99
$ws2.Range("E24").Value = 0
$ws2.Range("D26").Value = 457.92
$ws2.Range("E26").Value = 447.78
$ws2.Range("F26").Value = 0
$ws2.Range("C27").Value = 497.66
$ws2.Range("D27").Value = 8255.23
$ws2.Range("E27").Value = 1122.59
$ws2.Range("F27").Value = 0
$ws2.Range("C28").Value = 475.2
$ws2.Range("D28").Value = 950.4
$ws2.Range("E28").Value = 3503.69
$ws2.Range("F28").Value = 0
$ws2.Range("E29").Value = 366.83
$ws2.Range("F29").Value = 0
$ws2.Range("C30").Value = 9158.4
$ws2.Range("D30").Value = 10350.26
$ws2.Range("E30").Value = 681.79
$ws2.Range("F30").Value = 0
$ws2.Range("C31").Value = 222.32
$ws2.Range("D31").Value = 8733.540000000001
$ws2.Range("E31").Value = 9123.92
$ws2.Range("F31").Value = 0
$ws2.Range("C32").Value = 11275.94
$ws2.Range("D32").Value = 8872.17
$ws2.Range("E32").Value = 739.2
$ws2.Range("F32").Value = 0
$ws2.Range("C33").Value = 61.75
$ws2.Range("D33").Value = 0
$ws2.Range("C34").Value = 486.71
$ws2.Range("D34").Value = 10174.33
$ws2.Range("E34").Value = 10282.96
$ws2.Range("F34").Value = 0
$ws2.Range("C35").Value = 0
$ws2.Range("C36").Value = 518.4
$ws2.Range("D36").Value = 0
$ws2.Range("E40").Value = 746.3
$ws2.Range("F40").Value = 0
$ws2.Range("C41").Value = 7942.96
$ws2.Range("D41").Value = 8322.860000000001
$ws2.Range("E41").Value = 9587
$ws2.Range("F41").Value = 0
$ws2.Range("C42").Value = 0
$ws2.Range("D42").Value = 1428.84
$ws2.Range("E42").Value = 0
$ws2.Range("C44").Value = 660.24
$ws2.Range("D44").Value = 7011.36
$ws2.Range("E44").Value = 0
$ws2.Range("C46").Value = 295.63
$ws2.Range("D46").Value = 0
$ws2.Range("C47").Value = 1824.2
$ws2.Range("D47").Value = 1949.23
$ws2.Range("E47").Value = 879.99
$ws2.Range("F47").Value = 0
$ws2.Range("C48").Value = 1151.68
$ws2.Range("D48").Value = 4381.83
$ws2.Range("E48").Value = 5784.88
$ws2.Range("F48").Value = 0
$ws2.Range("C49").Value = 11132.19
$ws2.Range("D49").Value = 2781.14
$ws2.Range("E49").Value = 6664.04
$ws2.Range("F49").Value = 0
$ws2.Range("C50").Value = 1549.1
$ws2.Range("D50").Value = 71.62
$ws2.Range("E50").Value = 500.82
$ws2.Range("F50").Value = 0
$ws2.Range("C52").Value = 6923.38
$ws2.Range("D52").Value = 108.12
$ws2.Range("E52").Value = 4657.32
$ws2.Range("F52").Value = 0
$ws2.Range("C54").Value = 0
$ws2.Range("C55").Value = 0
$ws2.Range("D55").Value = 3995.75
$ws2.Range("E55").Value = 0
$ws2.Range("D56").Value = 869.53
$ws2.Range("E56").Value = 1035.56
$ws2.Range("F56").Value = 0
$ws2.Range("C57").Value = 142.2
$ws2.Range("D57").Value = 0
$ws2.Range("C58").Value = 0
$ws2.Range("D58").Value = 1556.27
$ws2.Range("E58").Value = 1034.72
$ws2.Range("F58").Value = 0
$ws2.Range("C60").Value = 81440.64
$ws2.Range("D60").Value = 91039.07000000001
$ws2.Range("E60").Value = 91874.38
$ws2.Range("F60").Value = 0

Write-Output "edits applied"